$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New culture/RGB rows appended to the table (rows 18-22)
$newData = @(
    @("arpitan", "(81, 114, 255)"),
    @("basque", "(41, 146, 17)"),
    @("walloon", "(131, 127, 37)"),
    @("picard", "(110, 120, 161)"),
    @("norman", "(243, 13, 13)")
)

$startRow = 18
for ($i = 0; $i -lt $newData.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $newData[$i][0]
    $ws.Cells.Item($r, 2).Value = $newData[$i][1]
}

# Update selection to match the new active cell after the edit session
$ws.Range("E15").Select() | Out-Null
